$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13140.909
$ws.Range("H23").Value = 13140.909
$ws.Range("H40").Value = 1484.2667
$ws.Range("I40").Value = 1469.7142
$ws.Range("J40").Value = 1488.6957
$ws.Range("K40").Value = 1469.7142
$ws.Range("L40").Value = 1488.6957
$ws.Range("M40").Value = -1294.7142
$ws.Range("N40").Value = -1838.6957
$ws.Range("H74").Value = 1002833.2
$ws.Range("I74").Value = 1002833.2
$ws.Range("K74").Value = 1002833.2
$ws.Range("M74").Value = -1001897.2
$ws.Range("H77").Value = 1002833.2
$ws.Range("I77").Value = 1002833.2
$ws.Range("K77").Value = 5014166
$ws.Range("M77").Value = -5009486
$ws.Range("H132").Value = 2615.2942
$ws.Range("I132").Value = 1841.25
$ws.Range("K132").Value = 5523.75
$ws.Range("M132").Value = -2993.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 956.7273
$ws.Range("I74").Value = 852.5
$ws.Range("K74").Value = 852.5
$ws.Range("M74").Value = 21.5
$ws.Range("H77").Value = 956.7273
$ws.Range("I77").Value = 852.5
$ws.Range("K77").Value = 4262.5
$ws.Range("M77").Value = 105.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2363.3157
$ws.Range("J86").Value = 2750
$ws.Range("L86").Value = 2750
$ws.Range("N86").Value = -4996
$ws.Range("H89").Value = 2363.3157
$ws.Range("J89").Value = 2750
$ws.Range("L89").Value = 13750
$ws.Range("N89").Value = -24982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 12500186
$ws.Range("I7").Value = 20000200
$ws.Range("K7").Value = 60000600
$ws.Range("M7").Value = -60000488
$ws.Range("H23").Value = 333346.34
$ws.Range("J23").Value = 999999
$ws.Range("L23").Value = 2999997
$ws.Range("N23").Value = -3000467
$ws.Range("H28").Value = 1065
$ws.Range("I28").Value = 1065
$ws.Range("K28").Value = 3195
$ws.Range("M28").Value = -2963
$ws.Range("H32").Value = 8712572
$ws.Range("J32").Value = 8712572
$ws.Range("L32").Value = 26137716
$ws.Range("N32").Value = -26138282
$ws.Range("H36").Value = 445
$ws.Range("I36").Value = 445
$ws.Range("K36").Value = 1335
$ws.Range("M36").Value = -1166
$ws.Range("H44").Value = 2166.6667
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 2250
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 6750
$ws.Range("M44").Value = -5602
$ws.Range("N44").Value = -7546
$ws.Range("H56").Value = 12180.583
$ws.Range("I56").Value = 12180.583
$ws.Range("K56").Value = 12180.583
$ws.Range("M56").Value = -11650.583
$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 75000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -75540
$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 75000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -76872
$ws.Range("H81").Value = 1006.6667
$ws.Range("J81").Value = 1010
$ws.Range("L81").Value = 3030
$ws.Range("N81").Value = -5276
$ws.Range("H84").Value = 1006.6667
$ws.Range("J84").Value = 1010
$ws.Range("L84").Value = 9090
$ws.Range("N84").Value = -20322
$ws.Range("H92").Value = 253.11111
$ws.Range("J92").Value = 218
$ws.Range("L92").Value = 654
$ws.Range("N92").Value = -3150
$ws.Range("H109").Value = 1218.4117
$ws.Range("J109").Value = 2399.5
$ws.Range("L109").Value = 7198.5
$ws.Range("N109").Value = -9278.5
$ws.Range("H110").Value = 5770.75
$ws.Range("I110").Value = 4351
$ws.Range("K110").Value = 13053
$ws.Range("M110").Value = -8963
$ws.Range("H112").Value = 2548.375
$ws.Range("I112").Value = 1937.4
$ws.Range("K112").Value = 5812.200000000001
$ws.Range("M112").Value = -4704.200000000001
$ws.Range("H120").Value = 7299.75
$ws.Range("I120").Value = 4099.5
$ws.Range("K120").Value = 12298.5
$ws.Range("M120").Value = -7460.5
$ws.Range("H131").Value = 1495.9623
$ws.Range("I131").Value = 1096.6666
$ws.Range("J131").Value = 1519.92
$ws.Range("K131").Value = 3289.9998
$ws.Range("L131").Value = 4559.76
$ws.Range("M131").Value = 1750.0002
$ws.Range("N131").Value = -14639.76
$ws.Range("H138").Value = 2425.25
$ws.Range("I138").Value = 1377.6
$ws.Range("K138").Value = 4132.799999999999
$ws.Range("M138").Value = 1007.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 22750126
$ws.Range("I113").Value = 41680300
$ws.Range("J113").Value = 33915.6
$ws.Range("K113").Value = 41680300
$ws.Range("L113").Value = 33915.6
$ws.Range("M113").Value = -41678130
$ws.Range("N113").Value = -38255.6
$ws.Range("H122").Value = 93490.37
$ws.Range("I122").Value = 2299.5
$ws.Range("K122").Value = 6898.5
$ws.Range("M122").Value = -4448.5
$ws.Range("H126").Value = 4222.25
$ws.Range("I126").Value = 4222.25
$ws.Range("K126").Value = 12666.75
$ws.Range("M126").Value = -10196.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16799
$ws.Range("I7").Value = 16799
$ws.Range("K7").Value = 16799
$ws.Range("M7").Value = -16687
$ws.Range("H126").Value = 16799
$ws.Range("I126").Value = 16799
$ws.Range("K126").Value = 50397
$ws.Range("M126").Value = -47927
